# Add data for 2024-11-17
# Updates CTA violent crime YTD figures across the citywide, by-neighborhood,
# and individual neighborhood sheets (2024 column plus a few historical revisions).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("F2").Value = 84
$ws.Range("B3").Value = 74
$ws.Range("H3").Value = 144
$ws.Range("B6").Value = 357
$ws.Range("C6").Value = 455
$ws.Range("E6").Value = 440
$ws.Range("F6").Value = 488
$ws.Range("K6").Value = 477
$ws.Range("B7").Value = 480
$ws.Range("C7").Value = 604
$ws.Range("E7").Value = 657
$ws.Range("F7").Value = 706
$ws.Range("H7").Value = 687
$ws.Range("K7").Value = 842

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("F6").Value = 21
$ws.Range("F7").Value = 47

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("F8").Value = 44
$ws.Range("K10").Value = 6
$ws.Range("E18").Value = 6
$ws.Range("B23").Value = 8
$ws.Range("H23").Value = 7
$ws.Range("K27").Value = 14
$ws.Range("K28").Value = 60
$ws.Range("F36").Value = 47
$ws.Range("B53").Value = 48
$ws.Range("F53").Value = 73
$ws.Range("K53").Value = 85
$ws.Range("C65").Value = 20
$ws.Range("F76").Value = 18
$ws.Range("E91").Value = 6
$ws.Range("B98").Value = 480
$ws.Range("C98").Value = 604
$ws.Range("E98").Value = 657
$ws.Range("F98").Value = 706
$ws.Range("H98").Value = 687
$ws.Range("K98").Value = 842

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("B6").Value = 30
$ws.Range("F6").Value = 53
$ws.Range("K6").Value = 44
$ws.Range("B7").Value = 48
$ws.Range("F7").Value = 73
$ws.Range("K7").Value = 85

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("E6").Value = 5
$ws.Range("E7").Value = 6

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("F5").Value = 15
$ws.Range("F6").Value = 18

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("C5").Value = 17
$ws.Range("C6").Value = 20

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("B3").Value = 1
$ws.Range("H3").Value = 2
$ws.Range("B7").Value = 8
$ws.Range("H7").Value = 7

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 6

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K5").Value = 3
$ws.Range("K6").Value = 6

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K4").Value = 10
$ws.Range("K5").Value = 14

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("F2").Value = 9
$ws.Range("F7").Value = 44
